$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plano de Ação")

# Fix small wording/casing issues on the existing row 40 (SPRINT 2D section)
$ws.Range("C40").Value = "Todos da equipe"
$ws.Range("H40").Value = "Criação da dashboard e linkar junto ao login "

# New requirement #1 (row 41): "Diagrama de solução"
$ws.Range("B41").Value = "Diagrama de solução "
$ws.Range("C41").Value = "Todos da equipe"
$ws.Range("D41").Value = "Essencial"
$ws.Range("E41").Value = 0.4
$ws.Range("E41").NumberFormat = $ws.Range("E40").NumberFormat
$ws.Range("F41").Value = 45219
$ws.Range("G41").Value = 45222
$ws.Range("H41").Value = "Terminar a criação do diagrama de solução"

# New requirement #2 (row 42): "slides da apresentação"
$ws.Range("B42").Value = "slides da apresentação"
$ws.Range("C42").Value = "Todos  da equipe"
$ws.Range("D42").Value = "Essencial"
$ws.Range("E42").Value = 0
$ws.Range("E42").NumberFormat = $ws.Range("E30").NumberFormat
$ws.Range("F42").Value = 45219
$ws.Range("G42").Value = 45222
$ws.Range("H42").Value = "Fazer os slides para a apresentação"

# Update the view state: scrolled frozen pane + active selection, matching the
# final on-screen state left after entering the new rows.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 36
$ws.Range("J42").Select()
